$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "Circuitos Elétricos 2"
$ws.Range("C3").Value = "Acionamentos"
$ws.Range("F3").Value = "Sistemas digitais"

# Row 4 (swap Acionamentos/Circuitos Elétricos 2 between C4 and D4)
$ws.Range("C4").Value = "Circuitos Elétricos 2"
$ws.Range("D4").Value = "Acionamentos"

# Row 6
$ws.Range("C6").Value = "Circuitos Elétricos 2"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "Sistemas digitais"
